$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the BTec logo images (headers) from image1.jpg -> image2.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

# Rename the Pearson logo images (footers) from image2.png -> image1.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -like "*PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}
